$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price text stays as text (matches original inline-string formatting)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "69.062.31"
$ws.Range("E2").Value = "  -1.73%  "

# Row 3
$ws.Range("D3").Value = "2.493.06"
$ws.Range("E3").Value = "  -1.13%  "

# Row 4
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.07%  "

# Row 5
$ws.Range("D5").Value = "568.61"
$ws.Range("E5").Value = "  -1.20%  "

# Row 6
$ws.Range("D6").Value = "164.31"
$ws.Range("E6").Value = "  -3.15%  "

# Row 7
$ws.Range("E7").Value = "  -0.08%  "

# Row 8
$ws.Range("D8").Value = "0.510"
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "2.495.74"
$ws.Range("E9").Value = "  -1.03%  "

# Row 10
$ws.Range("D10").Value = "0.158"
$ws.Range("E10").Value = "  -1.11%  "

# Row 11
$ws.Range("E11").Value = "  -0.51%  "

# Row 12
$ws.Range("D12").Value = "0.352"
$ws.Range("E12").Value = "  +2.52%  "

# Row 13
$ws.Range("D13").Value = "4.86"
$ws.Range("E13").Value = "  +0.79%  "

# Row 14
$ws.Range("D14").Value = "2.955.17"
$ws.Range("E14").Value = "  -0.96%  "

# Row 15
$ws.Range("D15").Value = "68.937.33"
$ws.Range("E15").Value = "  -1.75%  "

# Row 16
$ws.Range("D16").Value = "0.0000173"
$ws.Range("E16").Value = "  -2.99%  "

# Row 17
$ws.Range("D17").Value = "24.57"
$ws.Range("E17").Value = "  -1.93%  "

# Row 18
$ws.Range("D18").Value = "2.505.09"
$ws.Range("E18").Value = "  -0.61%  "

# Row 19
$ws.Range("D19").Value = "11.18"
$ws.Range("E19").Value = "  -3.10%  "

# Row 20
$ws.Range("D20").Value = "7.60"
$ws.Range("E20").Value = "  -0.05%  "

# Row 21
$ws.Range("D21").Value = "345.53"
$ws.Range("E21").Value = "  -2.82%  "

# Row 22
$ws.Range("D22").Value = "3.87"
$ws.Range("E22").Value = "  -2.15%  "

# Row 23
$ws.Range("D23").Value = "1.97"
$ws.Range("E23").Value = "  -1.03%  "

# Row 24
$ws.Range("E24").Value = "  +0.03%  "

# Row 25
$ws.Range("D25").Value = "69.91"
$ws.Range("E25").Value = "  +1.25%  "

# Row 26
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -5.26%  "

# Row 27
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.649.15"
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "8.77"
$ws.Range("E28").Value = "  -4.60%  "

# Row 29
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.08%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0875"
$ws.Range("E30").Value = "  -4.07%  "

# Row 31
$ws.Range("D31").Value = "7.70"
$ws.Range("E31").Value = "  -1.92%  "

# Row 32
$ws.Range("D32").Value = "455.74"
$ws.Range("E32").Value = "  -5.77%  "

# Row 33
$ws.Range("D33").Value = "1.21"
$ws.Range("E33").Value = "  -6.83%  "

# Row 34
$ws.Range("D34").Value = "1.71"
$ws.Range("E34").Value = "  -2.69%  "

# Row 35
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.11%  "

# Row 36
$ws.Range("E36").Value = "  -1.04%  "

# Row 37
$ws.Range("D37").Value = "155.50"
$ws.Range("E37").Value = "  -0.02%  "

# Row 38
$ws.Range("E38").Value = "  +0.61%  "

# Row 39
$ws.Range("D39").Value = "18.31"
$ws.Range("E39").Value = "  -1.58%  "

# Row 40
$ws.Range("E40").Value = "  -0.04%  "

# Row 41
$ws.Range("D41").Value = "0.315"
$ws.Range("E41").Value = "  -1.82%  "

# Row 42
$ws.Range("D42").Value = "4.63"
$ws.Range("E42").Value = "  -2.38%  "

# Row 43
$ws.Range("D43").Value = "1.58"
$ws.Range("E43").Value = "  -3.94%  "

# Row 44
$ws.Range("D44").Value = "38.07"
$ws.Range("E44").Value = "  -0.54%  "

# Row 45
$ws.Range("D45").Value = "1.11"
$ws.Range("E45").Value = "  -9.11%  "

# Row 46
$ws.Range("D46").Value = "2.18"
$ws.Range("E46").Value = "  -8.49%  "

# Row 47
$ws.Range("D47").Value = "140.88"
$ws.Range("E47").Value = "  -1.16%  "

# Row 48
$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D48").Value = "3.43"
$ws.Range("E48").Value = "  -2.65%  "

# Row 49
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "0.512"
$ws.Range("E49").Value = "  -2.71%  "

# Row 50
$ws.Range("D50").Value = "0.0728"
$ws.Range("E50").Value = "  -0.37%  "

# Row 51
$ws.Range("E51").Value = "  -4.22%  "
